$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 828.67645
$ws.Range("I15").Value = 828.67645
$ws.Range("K15").Value = 2486.02935
$ws.Range("M15").Value = -2317.02935
# Row 19
$ws.Range("H19").Value = 2623
$ws.Range("I19").Value = 2196.7144
$ws.Range("J19").Value = 2954.5557
$ws.Range("K19").Value = 2196.7144
$ws.Range("L19").Value = 2954.5557
$ws.Range("M19").Value = -2021.7144
$ws.Range("N19").Value = -3304.5557
# Row 33
$ws.Range("H33").Value = 677.5
$ws.Range("I33").Value = 175
$ws.Range("J33").Value = 2041.4286
$ws.Range("K33").Value = 175
$ws.Range("L33").Value = 2041.4286
$ws.Range("M33").Value = 54
$ws.Range("N33").Value = -2499.4286
# Row 34
$ws.Range("H34").Value = 14956.75
$ws.Range("I34").Value = 9109.166999999999
$ws.Range("J34").Value = 32499.5
$ws.Range("K34").Value = 9109.166999999999
$ws.Range("L34").Value = 32499.5
$ws.Range("M34").Value = -8906.166999999999
$ws.Range("N34").Value = -32905.5
# Row 36
$ws.Range("H36").Value = 14956.75
$ws.Range("I36").Value = 9109.166999999999
$ws.Range("J36").Value = 32499.5
$ws.Range("K36").Value = 9109.166999999999
$ws.Range("L36").Value = 32499.5
$ws.Range("M36").Value = -8394.166999999999
$ws.Range("N36").Value = -33929.5
# Row 40
$ws.Range("H40").Value = 5898.364
$ws.Range("J40").Value = 5349.5
$ws.Range("L40").Value = 5349.5
$ws.Range("N40").Value = -5699.5
# Row 43
$ws.Range("H43").Value = 1493.375
$ws.Range("I43").Value = 994
$ws.Range("K43").Value = 994
$ws.Range("M43").Value = -925
# Row 48
$ws.Range("H48").Value = 4829.8335
$ws.Range("J48").Value = 4829.8335
$ws.Range("L48").Value = 14489.5005
$ws.Range("N48").Value = -15073.5005
# Row 56
$ws.Range("H56").Value = 4829.8335
$ws.Range("J56").Value = 4829.8335
$ws.Range("L56").Value = 14489.5005
$ws.Range("N56").Value = -15557.5005
# Row 80
$ws.Range("H80").Value = 1005
$ws.Range("I80").Value = 845
$ws.Range("K80").Value = 2535
$ws.Range("M80").Value = -1537
# Row 83
$ws.Range("H83").Value = 1005
$ws.Range("I83").Value = 845
$ws.Range("K83").Value = 7605
$ws.Range("M83").Value = -2613
# Row 137
$ws.Range("H137").Value = 2684.4
$ws.Range("I137").Value = 2480.5
$ws.Range("K137").Value = 7441.5
$ws.Range("M137").Value = -4891.5
# Row 141
$ws.Range("H141").Value = 2466.8108
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6284.5293
$ws.Range("I32").Value = 824.7818
$ws.Range("K32").Value = 824.7818
$ws.Range("M32").Value = -537.7818
# Row 74
$ws.Range("H74").Value = 2094.6296
$ws.Range("I74").Value = 2113.6538
$ws.Range("K74").Value = 2113.6538
$ws.Range("M74").Value = -1239.6538
# Row 77
$ws.Range("H77").Value = 2094.6296
$ws.Range("I77").Value = 2113.6538
$ws.Range("K77").Value = 10568.269
$ws.Range("M77").Value = -6200.269
# Row 122
$ws.Range("H122").Value = 3539.8
$ws.Range("I122").Value = 3674.75
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 11024.25
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -8574.25
$ws.Range("N122").Value = -13900
# Row 132
$ws.Range("H132").Value = 2075.5
$ws.Range("I132").Value = 2118
$ws.Range("K132").Value = 6354
$ws.Range("M132").Value = -3824

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1959.7097
$ws.Range("I94").Value = 1622.75
$ws.Range("J94").Value = 2076.913
$ws.Range("K94").Value = 1622.75
$ws.Range("L94").Value = 2076.913
$ws.Range("M94").Value = -1171.75
$ws.Range("N94").Value = -2978.913
# Row 134
$ws.Range("H134").Value = 2737.9644
$ws.Range("I134").Value = 2737.9644
$ws.Range("K134").Value = 8213.893199999999
$ws.Range("M134").Value = -5678.893199999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 18
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
# Row 62
$ws.Range("H62").Value = 2000
$ws.Range("I62").Value = 2000
$ws.Range("K62").Value = 2000
$ws.Range("M62").Value = -1376
# Row 65
$ws.Range("H65").Value = 2000
$ws.Range("I65").Value = 2000
$ws.Range("K65").Value = 10000
$ws.Range("M65").Value = -6880
# Row 69
$ws.Range("H69").Value = 49999
$ws.Range("I69").Value = 49999
$ws.Range("K69").Value = 49999
$ws.Range("M69").Value = -49250
# Row 72
$ws.Range("H72").Value = 49999
$ws.Range("I72").Value = 49999
$ws.Range("K72").Value = 149997
$ws.Range("M72").Value = -146253
# Row 99
$ws.Range("H99").Value = 10007.637
$ws.Range("I99").Value = 7599
$ws.Range("K99").Value = 7599
$ws.Range("M99").Value = -6101
# Row 122
$ws.Range("H122").Value = 4999.5
$ws.Range("J122").Value = 4999
$ws.Range("L122").Value = 14997
$ws.Range("N122").Value = -19897
# Row 126
$ws.Range("H126").Value = 10007.637
$ws.Range("I126").Value = 7599
$ws.Range("K126").Value = 22797
$ws.Range("M126").Value = -20327
# Row 138
$ws.Range("H138").Value = 89997.5
$ws.Range("J138").Value = 89997.5
$ws.Range("L138").Value = 89997.5
$ws.Range("N138").Value = -100277.5
# Row 139
$ws.Range("H139").Value = 42531.953
$ws.Range("I139").Value = 31983.834
$ws.Range("J139").Value = 89998.5
$ws.Range("K139").Value = 31983.834
$ws.Range("L139").Value = 89998.5
$ws.Range("M139").Value = -26843.834
$ws.Range("N139").Value = -100278.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 465.63635
$ws.Range("I5").Value = 481.7
$ws.Range("K5").Value = 1445.1
$ws.Range("M5").Value = -1333.1
# Row 26
$ws.Range("H26").Value = 1934.2
$ws.Range("I26").Value = 1934.2
$ws.Range("K26").Value = 5802.6
$ws.Range("M26").Value = -5514.6
# Row 107
$ws.Range("H107").Value = 978
$ws.Range("I107").Value = 666.3333
$ws.Range("K107").Value = 1998.9999
$ws.Range("M107").Value = -78.99990000000003
# Row 122
$ws.Range("H122").Value = 2510.2354
$ws.Range("J122").Value = 3572.7273
$ws.Range("L122").Value = 32154.5457
$ws.Range("N122").Value = -37054.5457
# Row 135
$ws.Range("H135").Value = 465.63635
$ws.Range("I135").Value = 481.7
$ws.Range("K135").Value = 4335.3
$ws.Range("M135").Value = -1800.3

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 10000
$ws.Range("I5").Value = 10000
$ws.Range("K5").Value = 10000
$ws.Range("M5").Value = -9888
# Row 22
$ws.Range("H22").Value = 2500
$ws.Range("I22").Value = 2500
$ws.Range("K22").Value = 2500
$ws.Range("M22").Value = -1971
# Row 80
$ws.Range("H80").Value = 3428.4443
$ws.Range("I80").Value = 3582.4285
$ws.Range("K80").Value = 3582.4285
$ws.Range("M80").Value = -2584.4285
# Row 83
$ws.Range("H83").Value = 3428.4443
$ws.Range("I83").Value = 3582.4285
$ws.Range("K83").Value = 17912.1425
$ws.Range("M83").Value = -12920.1425
# Row 126
$ws.Range("H126").Value = 2011.4166
$ws.Range("I126").Value = 1347.5
$ws.Range("J126").Value = 2144.2
$ws.Range("K126").Value = 4042.5
$ws.Range("L126").Value = 6432.599999999999
$ws.Range("M126").Value = -1572.5
$ws.Range("N126").Value = -11372.6
# Row 132
$ws.Range("H132").Value = 3872.7593
$ws.Range("I132").Value = 4048.535
$ws.Range("K132").Value = 12145.605
$ws.Range("M132").Value = -9615.605

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1023.9583
$ws.Range("I16").Value = 1030.8422
$ws.Range("J16").Value = 997.8
$ws.Range("K16").Value = 1030.8422
$ws.Range("L16").Value = 997.8
$ws.Range("M16").Value = -860.8422
$ws.Range("N16").Value = -1337.8
# Row 40
$ws.Range("H40").Value = 1852.3334
$ws.Range("I40").Value = 1824.8572
$ws.Range("J40").Value = 1948.5
$ws.Range("K40").Value = 1824.8572
$ws.Range("L40").Value = 1948.5
$ws.Range("M40").Value = -1688.8572
$ws.Range("N40").Value = -2220.5
# Row 68
$ws.Range("H68").Value = 2200
$ws.Range("I68").Value = 1639.4
$ws.Range("K68").Value = 1639.4
$ws.Range("M68").Value = -890.4000000000001
# Row 71
$ws.Range("H71").Value = 2200
$ws.Range("I71").Value = 1639.4
$ws.Range("K71").Value = 8197
$ws.Range("M71").Value = -4453
# Row 132
$ws.Range("H132").Value = 2423.074
$ws.Range("I132").Value = 2150.0527
$ws.Range("K132").Value = 6450.158100000001
$ws.Range("M132").Value = -3920.158100000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 203299.8
$ws.Range("I5").Value = 5999
$ws.Range("K5").Value = 5999
$ws.Range("M5").Value = -5887
# Row 45
$ws.Range("H45").Value = 35413.7
$ws.Range("I45").Value = 44828.5
$ws.Range("J45").Value = 29137.166
$ws.Range("K45").Value = 44828.5
$ws.Range("L45").Value = 29137.166
$ws.Range("M45").Value = -44337.5
$ws.Range("N45").Value = -30119.166
# Row 81
$ws.Range("H81").Value = 2998.25
$ws.Range("I81").Value = 2998.25
$ws.Range("K81").Value = 5996.5
$ws.Range("M81").Value = -4935.5
# Row 84
$ws.Range("H84").Value = 2998.25
$ws.Range("I84").Value = 2998.25
$ws.Range("K84").Value = 29982.5
$ws.Range("M84").Value = -24678.5
# Row 132
$ws.Range("H132").Value = 4217.9355
$ws.Range("I132").Value = 3502.2693
$ws.Range("J132").Value = 7939.4
$ws.Range("K132").Value = 10506.8079
$ws.Range("L132").Value = 23818.2
$ws.Range("M132").Value = -7976.8079
$ws.Range("N132").Value = -28878.2
# Row 136
$ws.Range("H136").Value = 1274.1428
$ws.Range("I136").Value = 1100.4634
$ws.Range("K136").Value = 3301.3902
$ws.Range("M136").Value = -751.3902000000003
